$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Remove the "License Information" (Heading2) paragraph entirely.
# ------------------------------------------------------------------
$rng = $d.Content
$rng.Find.ClearFormatting()
$found = $rng.Find.Execute("License Information", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $rng.Expand(4)
    $rng.Delete()
}

# ------------------------------------------------------------------
# 2. Rewrite the licensing paragraph that used to read:
#    "Translation Questions (unfoldingWord) is based on: unfoldingWord(r)
#     Translation Questions, unfoldingWord, 2022, which is licensed under
#     a CC BY-SA 4.0 license."
# ------------------------------------------------------------------
$rng2 = $d.Content
$rng2.Find.ClearFormatting()
$found2 = $rng2.Find.Execute("Translation Questions (unfoldingWord) is based on", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found2) {
    $target = $rng2.Duplicate
    $target.Expand(4)
    # Trim the trailing paragraph mark so we only touch the visible text,
    # leaving the (already present) leading and trailing empty runs intact.
    $target.MoveEnd(1, -1)
    $target.Delete()

    $ins = $d.Range($target.Start, $target.Start)

    $ins.InsertAfter("unfoldingWord® Translation Questions")
    $ins.LanguageID = "en_US"
    $ins.LanguageIDOther = "en_US"
    $ins.Font.Bold = $true
    $ins.Collapse(0)

    $ins.InsertAfter(" © 2022 unfoldingWord. Released under CC BY-SA 4.0 license. ")
    $ins.LanguageID = "en_US"
    $ins.LanguageIDOther = "en_US"
    $ins.Font.Bold = $false
    $ins.Collapse(0)

    $ins.InsertAfter("unfoldingWord® Translation Questions")
    $ins.LanguageID = "en_US"
    $ins.LanguageIDOther = "en_US"
    $ins.Font.Bold = $false
    $ins.Collapse(0)

    $ins.InsertAfter(" has been adapted in the following languages: Tok Pisin, Arabic (عربي), French (Français), Hindi (हिंदी), Indonesian (Bahasa Indonesia), Portuguese (Português), Russian (Русский), Spanish (Español), Swahili (Kiswahili), and Simplified Chinese (简体中文) from ")
    $ins.LanguageID = "en_US"
    $ins.LanguageIDOther = "en_US"
    $ins.Font.Bold = $false
    $ins.Collapse(0)

    $ins.InsertAfter("unfoldingWord® Translation Questions")
    $ins.LanguageID = "en_US"
    $ins.LanguageIDOther = "en_US"
    $ins.Font.Bold = $false
    $ins.Collapse(0)

    $ins.InsertAfter(" © 2022 unfoldingWord. Released under CC BY-SA 4.0 license by Mission Mutual")
    $ins.LanguageID = "en_US"
    $ins.LanguageIDOther = "en_US"
    $ins.Font.Bold = $false
}

# ------------------------------------------------------------------
# 3. Remove the "This PDF version is provided under the same license."
#    paragraph entirely.
# ------------------------------------------------------------------
$rng3 = $d.Content
$rng3.Find.ClearFormatting()
$found3 = $rng3.Find.Execute("This PDF version is provided under the same license.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found3) {
    $rng3.Expand(4)
    $rng3.Delete()
}

Write-Host "Done"
